$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.280.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.68%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.122.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.49%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.97%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.542"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.121.08"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.52%  "
$ws.Range("E10").Value = "  -5.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.95"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.456"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.33"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000241"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.634.70"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.59%  "
$ws.Range("E16").Value = "  -1.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.28"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.189.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.117.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "481.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.718"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.14%  "
$ws.Range("E28").Value = "  -0.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.02%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.35%  "
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("E33").Value = "  -6.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.51"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0847"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.24%  "
$ws.Range("E36").Value = "  -2.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.07"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.31"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -9.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.26"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "51.08"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.14%  "
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "453.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.68%  "
$ws.Range("B42").Value = "Cosmos"
$ws.Range("C42").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.293"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0367"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.34%  "
$ws.Range("E45").Value = "  -0.55%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.851.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.42%  "
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.57"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.01%  "
